$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 361; this pushes the existing data in
# rows 361:469 down to rows 362:470 and extends the sheet dimension to
# A1:R470 automatically.
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new weekly record.
$ws.Range("A361").Value = 9
$ws.Range("B361").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C361").Value = "Metropolitana"
$ws.Range("D361").Value = 44841
$ws.Range("E361").Value = 13
$ws.Range("F361").Value = 100112039
$ws.Range("G361").Value = "Ciboulette"
$ws.Range("H361").Value = "Sin especificar"
$ws.Range("I361").Value = "Primera"
$ws.Range("J361").Value = 250
$ws.Range("K361").Value = 800
$ws.Range("L361").Value = 1000
$ws.Range("M361").Value = 900
$ws.Range("N361").Value = "`$/docena de atados"
$ws.Range("O361").Value = "Región Metropolitana"
$ws.Range("P361").Value = 300
$ws.Range("Q361").Value = 3
$ws.Range("R361").Value = "Hortaliza"
